# Add a new vim cheat-sheet entry: "run bash cmd in vim"
# Appends two rows (52, 53) to the end of the "Command Mode" section on Sheet1:
#   row 52: "Run command"        | ":! Echo…"
#   row 53: "Search and replace" | ":s/foo/bar/g`n:%s/foo/bar/g (do all)"
# Also widens column C and gives the new multi-line cell centered/wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 52: plain entry, matches style of existing "Command Mode" rows ---
$ws.Range("B52").Value = "Run command"
$ws.Range("C52").Value = ":! Echo…"

# B52 picks up the same yellow-highlight look used elsewhere in this block
# (e.g. B20) so it lands on the existing fill style instead of a new one.
$ws.Range("B52").Interior.Color = 65535

# --- Row 53: new command with a two-line value in column C ---
$ws.Range("B53").Value = "Search and replace"
$ws.Range("C53").Value = ":s/foo/bar/g`n:%s/foo/bar/g (do all)"

# Give the wrapped two-line cell a taller row and centered/wrapped alignment.
$ws.Range("C53").HorizontalAlignment = -4108
$ws.Range("C53").WrapText = $true
$ws.Rows.Item(53).RowHeight = 30

# Column C needs to be wider now that it holds a two-line value.
$ws.Columns.Item(3).ColumnWidth = 25.1666667

# Reflect the new bottom-of-sheet selection, matching where the author
# was last working in the sheet.
$ws.Range("D52").Select()

Write-Output "done"
